# Revise the TableTemplate worksheet:
#  - Insert a new "Is Foreign Key" column between "Part of Primary Key" and "Data Type"
#  - Rename "Part of Key" -> "Part of Primary Key"
#  - Rename "Default" -> "Default value"
#  - Translate SQL data types to C#/EF Core data types (nvarchar -> string, bit -> bool)
#  - Populate the new "Is Foreign Key" column with "no" for every data row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E; this shifts the old Data Type/Allow Null/Default
# columns (E,F,G) one slot to the right (F,G,H), preserving their contents,
# styles and the "Person entity" value that lived in F3.
$ws.Columns("E:E").Insert()

# --- Header row (row 6) ---
$ws.Range("D6").Value = "Part of Primary Key"
$ws.Range("E6").Value = "Is Foreign Key"
$ws.Range("H6").Value = "Default value"

# --- New "Is Foreign Key" column values (all "no") ---
# E8 inherits D8's shaded style from the column Insert; reset it back to the
# workbook's default (unstyled) cell style to match the rest of the column.
$ws.Range("E8").Style = "Normal"
$ws.Range("E8").Value = "no"
$ws.Range("E9").Value = "no"
$ws.Range("E10").Value = "no"
$ws.Range("E11").Value = "no"
$ws.Range("E12").Value = "no"

# --- Translate SQL data types in the (shifted) "Data Type" column F ---
$ws.Range("F8").Value = "int"
$ws.Range("F9").Value = "string"
$ws.Range("F10").Value = "string"
$ws.Range("F11").Value = "int"
$ws.Range("F12").Value = "bool"

$ws.Range("G19").Select()
